$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "305.67") but must stay as literal
# text to match the source data (thousand-dot formatted prices, some of
# which - e.g. "44.047.52" - are not even valid numbers). Force text format,
# assign, then restore the Normal style so no stray formatting is left behind.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "44.047.52"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.237.80"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "305.67"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "95.63"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.524"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.78"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.22"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.578.10"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.238.71"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.822"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "13.57"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "44.000.99"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "64.75"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "236.46"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.89"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "37.44"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.13"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.94"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "20.03"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "154.97"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0806"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.32"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "15.06"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.38"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.82"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.741.25"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "84.73"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "100.10"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "4.94"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "14.70"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "68.98"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "8.08"
$c.Style = "Normal"

# Coin names / links / volume percentages are never numeric-looking, so a plain
# assignment is enough to keep them as text.
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  -4.31%  "
$ws.Range("E6").Value = "  -6.64%  "
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -5.63%  "
$ws.Range("E10").Value = "  -7.24%  "
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("E12").Value = "  -5.63%  "
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("E16").Value = "  -4.61%  "
$ws.Range("E17").Value = "  -6.28%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("E20").Value = "  -8.13%  "
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E24").Value = "  -7.50%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E25").Value = "  -7.75%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("E33").Value = "  -5.46%  "
$ws.Range("E34").Value = "  +8.93%  "
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("E36").Value = "  -5.67%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("E38").Value = "  -10.27%  "
$ws.Range("E39").Value = "  -10.24%  "
$ws.Range("E40").Value = "  -9.34%  "
$ws.Range("E41").Value = "  -9.42%  "
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("E46").Value = "  -6.40%  "
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("E50").Value = "  -9.57%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E51").Value = "  -4.25%  "
